$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 110, shifting the existing rows
# 110-167 down to 111-168 (keeps formatting/styles of the row above).
$ws.Range("A110").EntireRow.Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A110").Value = 5
$ws.Range("B110").Value = "Macroferia Regional de Talca"
$ws.Range("C110").Value = "Maule"
$ws.Range("D110").Value = 44452
$ws.Range("E110").Value = 7
$ws.Range("F110").Value = 100114014
$ws.Range("G110").Value = "Betarraga"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 4000
$ws.Range("K110").Value = 700
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 700
$ws.Range("N110").Value = '$/paquete 5 unidades'
$ws.Range("O110").Value = "Región del Maule"
$ws.Range("P110").Value = 140
$ws.Range("Q110").Value = 5
$ws.Range("R110").Value = "Hortaliza"
